$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph that ends the document's body text - the one ending in
# "...I think there is enough stuff here to figure out how to intermingle it
# all up and get cool results. " - which is immediately followed by a single
# trailing blank paragraph before the section break.
# ---------------------------------------------------------------------------
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*intermingle it all up*get cool results*") {
        $anchorIndex = $i
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not find the anchor paragraph (...get cool results.)"
}

# The paragraph right after the anchor is the document's final (blank) paragraph.
$trailingIndex = $anchorIndex + 1
$trailingRange = $d.Paragraphs.Item($trailingIndex).Range

# Remember where the new paragraphs will start (this index, and the next four,
# will hold: a blank paragraph, then the four new sentences/paragraphs). The
# original trailing blank paragraph gets pushed down below all of them.
$newStart = $trailingIndex

# Create 5 new empty paragraphs immediately before the existing trailing blank
# paragraph. Doing all the InsertParagraphBefore calls first (before filling
# in any text) keeps paragraph indices stable and predictable.
for ($k = 0; $k -lt 5; $k++) {
    $trailingRange.InsertParagraphBefore()
}

# $newStart       -> stays blank (matches the blank separator paragraph)
# $newStart + 1   -> "After I 've implemented ... propositions."
# $newStart + 2   -> "The first and major problem ... 0s."
# $newStart + 3   -> "It is needed to say ... and so on."
# $newStart + 4   -> "For that purpose, ... algorithm."
# $newStart + 5   -> original trailing blank paragraph (unchanged)

$d.Paragraphs.Item($newStart + 1).Range.InsertAfter("After I " + [char]0x2018 + "ve implemented all the aforementioned, I found one serious flaw of the above propositions.")

$d.Paragraphs.Item($newStart + 2).Range.InsertAfter("The first and major problem is that we specify only one value in our polynomial that is 0s.")

$d.Paragraphs.Item($newStart + 3).Range.InsertAfter("It is needed to say to the function that for example the letter " + [char]0x201C + "A" + [char]0x201D + " should give us results close to 1, letter " + [char]0x201C + "B" + [char]0x201D + " close to 2 and so on.")

$d.Paragraphs.Item($newStart + 4).Range.InsertAfter("For that purpose, one may use Lagrange" + [char]0x2019 + "s method for interpolation or Newton" + [char]0x2019 + "s Divided differences algorithm.")

Write-Output "Inserted new paragraphs after index $anchorIndex; document now has $($d.Paragraphs.Count) paragraphs."
